$d = $word.ActiveDocument

# Locate the start paragraph (contains the SaaS description intro) and the end
# paragraph (ends with the "...para la asignacion de trabajo." sentence) by
# scanning for unique anchor substrings, since paragraph indices can shift.
$startIdx = -1
$endIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    $t = $p.Range.Text
    if ($startIdx -eq -1 -and $t -like "*tipo de empresas clientes*Software as a Service*") {
        $startIdx = $i
    }
    if ($t -like "*Para la asignaci*n de carga de trabajo para los t*cnicos el sistema permitir**") {
        $endIdx = $i
    }
}

if ($startIdx -eq -1 -or $endIdx -eq -1) {
    throw "Could not locate target paragraphs (start=$startIdx end=$endIdx)"
}

$pStart = $d.Paragraphs($startIdx)
$pEnd = $d.Paragraphs($endIdx)

$rangeStart = $pStart.Range.Start
$rangeEnd = $pEnd.Range.End - 1

$target = $d.Range($rangeStart, $rangeEnd)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="709"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Según la descripción del problema y el tipo de empresas clientes, el proyecto tendrá las características de un “Software as a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Service</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SaaS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>), lo que significa que se proveerá de sistema como un servicio externo a la empresa donde esta organización  solo se preocupara de consumir el sistema, en ningún caso del desarrollo, mantención o reparación de este. La empresa cliente deberá ser responsable de poseer dispositivos que se conecten vía internet al sistema, además de la configuración y administración básica del sistema para que se acomode y personalice a las características de la empresa, por ejemplo, logo y nombre corporativo, usuarios, definición de contrato, normativas internas, etc.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t>Una característica importante es que el sistema será evolutivo (modular) y crecerá según demanda de nuevas opciones y mejoras que el o los clientes requieran.</w:t></w:r><w:r><w:br/><w:t>El proyecto es un sistema para el control de flujos de procesos de un servicio técnico electrónico. Para este procedimiento el sistema deberá permitir registrar los artículos que ingresan al servicio y administrar de manera eficiente la carga de trabajo para los técnicos.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:br/></w:r><w:r><w:tab/></w:r><w:r><w:t>En el registro de artículos que ingresan al servicio técnico deberá ser un mantenedor de órdenes de trabajo, además debe contar con un flujo de trabajo manejado por estados, los cuales indicarán la situación actual dentro del servicio.</w:t></w:r><w:r><w:br/><w:t>Para la administración en la carga de trabajo de los técnicos, el sistema deberá permitir realizar una asignación de tareas a los técnicos dependiendo de su especialidad y disponibilidad de tiempo al momento de realizar la carga. El sistema será capaz de balancear la carga de trabajo hacia los técnicos disponibles según especialidad, el tiempo estimado a reparar el artículo tiene un tope preestablecido o configurable e internamente se manejara la duración en que se desarrollaron las tareas para una mayor optimización en el uso de recursos y análisis estadísticos de desempeño.</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Para el registro de órdenes de trabajo el sistema debe permitir registrar todos los datos necesarios para identificar el artículo, además debe permitir realizar el flujo de trabajo de manera ordenada, y dependiente de los diferentes perfiles de usuarios con las autorizaciones correspondientes. Para el cambio de estados será responsabilidad de cada actor dentro del flujo de trabajo cumplir la labor de actualizar el estado correspondiente del artículo. Se restringirán los accesos definiendo los perfiles correspondientes que serán usados en el sistema y los estados también se regirán por el mismo modelo se asociaran a los perfiles. Para lograr el ajuste de cada perfil, el sistema debe permitir la configuración de perfiles con los permisos correspondientes. El </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>sistema debe permitir agregar y modificar los estados actuales.</w:t></w:r><w:r><w:br/><w:t>Para la asignación de carga de trabajo para los técnicos el sistema permitirá asignar diferentes especialidades y dependiendo de sus habilidades se diferenciaran en niveles para la asignación de trabajo.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$target.InsertXML($xml)

Write-Host "Replaced paragraphs $startIdx..$endIdx"
